# Atualizacao de bases das ligas, do dia: 24-02-2024 as 12:40
#
# The source data feed re-sorted a handful of fixtures that share the same
# match date. Columns A (id), C (Div), D (Div Original Name) and E (Date)
# stay put, while the per-match payload in columns B and F:AC is rotated
# among the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows($sheet, $rows) {
    # Snapshot B:AC for every row in the block before writing anything back,
    # otherwise later writes would clobber values still needed for earlier ones.
    $snapshots = @{}
    foreach ($r in $rows) {
        $snapshots[$r] = $sheet.Range("B$r`:AC$r").Value2
    }

    $count = $rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $rows[$i]
        # destination row i receives the data that used to live in the
        # previous row of the block (wrapping around to the last row).
        $srcRow = $rows[($i - 1 + $count) % $count]
        $sheet.Range("B$destRow`:AC$destRow").Value2 = $snapshots[$srcRow]
    }
}

Rotate-Rows $ws @(21, 22, 23)
Rotate-Rows $ws @(44, 45, 46)
Rotate-Rows $ws @(104, 105)
Rotate-Rows $ws @(119, 120)
